$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.148.23'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.487.53'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.81'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.29'
$ws.Range("E6").Value = '  +2.90%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("D9").Value = '2.487.40'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("E14").Value = '  +1.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.32'
$ws.Range("E15").Value = '  -2.38%  '
$ws.Range("D16").Value = '67.056.57'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("D18").Value = '2.572.13'
$ws.Range("E18").Value = '  +2.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.00'
$ws.Range("E19").Value = '  -5.25%  '
$ws.Range("E20").Value = '  -5.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.34'
$ws.Range("E21").Value = '  -3.22%  '
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("B24").Value = 'NEARProtocol'
$ws.Range("C24").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.23'
$ws.Range("E24").Value = '  -4.17%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '68.41'
$ws.Range("E25").Value = '  -2.98%  '
$ws.Range("E26").Value = '  -2.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.29'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("E30").Value = '  -3.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '510.43'
$ws.Range("E31").Value = '  +2.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.73'
$ws.Range("E32").Value = '  -3.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  -2.88%  '
$ws.Range("E34").Value = '  -3.95%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.89'
$ws.Range("E36").Value = '  +0.96%  '
$ws.Range("E37").Value = '  -7.53%  '
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.24'
$ws.Range("E39").Value = '  -4.24%  '
$ws.Range("E40").Value = '  -5.63%  '
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("E44").Value = '  -2.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.36'
$ws.Range("E45").Value = '  -4.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.79'
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.65'
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("E48").Value = '  -4.37%  '
$ws.Range("E49").Value = '  -4.82%  '
$ws.Range("E50").Value = '  -4.88%  '
$ws.Range("E51").Value = '  -0.76%  '
